# Update dashboard with new charts/data
# - Rename the "06.09.25" tab to "25-09-2025" (new report date).
# - Zero out Rufus Kirimi Ntongai's (row 16) pre-payment balance, which
#   cascades through the row's SUM formulas (G16, K16, L16) and the
#   Grand Total row (row 19: B19, G19, K19, L19) automatically.
# - Move the viewport/selection to where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first worksheet to reflect the new balances date.
$ws.Name = "25-09-2025"

# Make sure it's the active sheet before touching selection/scroll state.
$ws.Activate()

# Rufus Kirimi Ntongai's pre-payment (column B) balance moved to 0.
$ws.Range("B16").Value = 0

# Scroll the view down (topLeftCell ~ A10) and leave the selection on B35,
# matching where editing continued.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B35").Select()

$wb.Save()
